$wb = $excel.ActiveWorkbook

# Column F ("想去人数") updates that apply identically to both the
# "展览" and "全部类型" worksheets (rows 2-11, row 7 and row 12 unchanged).
$updates = @{
    2  = 2834
    3  = 730
    4  = 96
    5  = 6680
    6  = 1573
    8  = 27
    9  = 44
    10 = 102
    11 = 17
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
